# Update "想去人数" (want-to-go count) figures in F4, F6, F8
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1446
    $ws.Range("F6").Value = 27
    $ws.Range("F8").Value = 39
}
